$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '47.631.21'
$ws.Range('E2').Value = '  +0.15%  '
$ws.Range('D3').Value = '2.495.05'
$ws.Range('E3').Value = '  +0.26%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '321.73'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '108.96'
$ws.Range('E6').Value = '  +3.77%  '
$ws.Range('E7').Value = '  -0.57%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('E9').Value = '  -0.03%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '39.33'
$ws.Range('E10').Value = '  +2.56%  '
$ws.Range('E11').Value = '  -0.58%  '
$ws.Range('E12').Value = '  +0.51%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '18.67'
$ws.Range('E13').Value = '  +2.01%  '
$ws.Range('E14').Value = '  +0.61%  '
$ws.Range('D15').Value = '2.885.73'
$ws.Range('E15').Value = '  +0.30%  '
$ws.Range('D16').Value = '2.501.04'
$ws.Range('E16').Value = '  +0.22%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.848'
$ws.Range('E17').Value = '  +0.16%  '
$ws.Range('D18').Value = '47.487.61'
$ws.Range('E18').Value = '  +0.04%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '13.35'
$ws.Range('E19').Value = '  +4.71%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '6.64'
$ws.Range('D21').Value = '0.0₃0940'
$ws.Range('E21').Value = '  +0.34%  '
$ws.Range('E22').Value = '  +14.83%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '70.65'
$ws.Range('E23').Value = '  -0.04%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '246.95'
$ws.Range('E24').Value = '  -1.68%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.55'
$ws.Range('E25').Value = '  -1.18%  '
$ws.Range('E26').Value = '  +0.04%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '25.77'
$ws.Range('E27').Value = '  -1.45%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.00'
$ws.Range('E28').Value = '  -0.38%  '
$ws.Range('B29').Value = 'Toncoin'
$ws.Range('C29').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.20'
$ws.Range('E29').Value = '  -0.05%  '
$ws.Range('B30').Value = 'Kaspa'
$ws.Range('C30').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.140'
$ws.Range('E30').Value = '  +4.04%  '
$ws.Range('B31').Value = 'InjectiveProtocol'
$ws.Range('C31').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '34.67'
$ws.Range('E31').Value = '  -1.36%  '
$ws.Range('E32').Value = '  +1.06%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '20.19'
$ws.Range('E33').Value = '  +2.09%  '
$ws.Range('E34').Value = '  -0.61%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.0788'
$ws.Range('E35').Value = '  +0.61%  '
$ws.Range('E36').Value = '  +0.16%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '4.72'
$ws.Range('E37').Value = '  +1.70%  '
$ws.Range('E38').Value = '  -0.49%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.95'
$ws.Range('E39').Value = '  -1.29%  '
$ws.Range('E40').Value = '  +0.21%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '22.18'
$ws.Range('E41').Value = '  +3.76%  '
$ws.Range('E42').Value = '  -2.15%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '119.03'
$ws.Range('E43').Value = '  -2.89%  '
$ws.Range('E44').Value = '  -0.21%  '
$ws.Range('D45').Value = '1.997.04'
$ws.Range('E45').Value = '  +1.53%  '
$ws.Range('E46').Value = '  +1.76%  '
$ws.Range('E47').Value = '  -2.91%  '
$ws.Range('E48').Value = '  -0.16%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '9.06'
$ws.Range('E49').Value = '  -1.59%  '
$ws.Range('E50').Value = '  -0.75%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '56.66'
$ws.Range('E51').Value = '  +3.37%  '
